# Add three new country sheets (Russia, Finland, Hungary) to the workbook,
# following the same "country market" template used by the existing sheets.
#
# The three sheets are copied from existing sheets that already have the
# exact same layout/column-widths/row-heights needed:
#   - Russia / Finland  <- copied from "Croatia"  (narrower column B)
#   - Hungary           <- copied from "Turkey"   (wider column B)
#
# A throw-away sheet is created (and deleted at the end) purely so the
# workbook's internal sheetId counter advances the same way it did in the
# original authoring session (the three new sheets end up with sheetId
# 22/23/24, leaving 21 "spent").

$wb = $excel.ActiveWorkbook

$scratchSource = $wb.Worksheets.Item("Croatia")
$endOfBook = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratchSource.Copy($null, $endOfBook)
$scratch = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---- Russia --------------------------------------------------------------
$russiaSource = $wb.Worksheets.Item("Croatia")
$endOfBook = $wb.Worksheets.Item($wb.Worksheets.Count)
$russiaSource.Copy($null, $endOfBook)
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T3313"
$russia.Range("B2").Value = "Russia Market"

# ---- Finland --------------------------------------------------------------
$finlandSource = $wb.Worksheets.Item("Croatia")
$endOfBook = $wb.Worksheets.Item($wb.Worksheets.Count)
$finlandSource.Copy($null, $endOfBook)
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2890"
$finland.Range("B2").Value = "Finland Market"

# ---- Hungary ---------------------------------------------------------------
$hungarySource = $wb.Worksheets.Item("Turkey")
$endOfBook = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungarySource.Copy($null, $endOfBook)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T2999/T2982"
$hungary.Range("B2").Value = "Hungary Market"

# Drop the scratch sheet now that the sheetId counter has advanced past it.
$scratch.Delete()

# Re-fetch fresh worksheet references by name (the delete above can leave
# previously-held references to other sheets with a stale "active" state)
# and restore each sheet's own selection/view, matching the final workbook.
$russia = $wb.Worksheets.Item("Russia")
$russia.Activate()
$russia.Range("A1:D12").Select()

$finland = $wb.Worksheets.Item("Finland")
$finland.Activate()
$finland.Range("A1:D12").Select()

$hungary = $wb.Worksheets.Item("Hungary")
$hungary.Activate()
$hungary.Range("I14").Select()
